$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Reorder the HKL bracket labels in the header row (row 2, C:I) ---
# Old order: [1,1,0] [2,0,0] [2,1,1] [2,2,0] [2,2,2] [3,1,0] [3,2,1]   (J2 [4,0,0] unchanged)
# New order: [3,2,1] [1,1,0] [3,1,0] [2,2,2] [2,0,0] [2,2,0] [2,1,1]
$ws.Range("C2").Value = "[3, 2, 1]"
$ws.Range("D2").Value = "[1, 1, 0]"
$ws.Range("E2").Value = "[3, 1, 0]"
$ws.Range("F2").Value = "[2, 2, 2]"
$ws.Range("G2").Value = "[2, 0, 0]"
$ws.Range("H2").Value = "[2, 2, 0]"
$ws.Range("I2").Value = "[2, 1, 1]"

# --- Append 4 new data rows (20-23) for the Holden scheme, copying the
#     formatting pattern already used by the existing data rows ---
$newRows = @(
    @{ Row = 20; A = 18; B = "HexGrid-90degTilt2.5degRes" },
    @{ Row = 21; A = 19; B = "HexGrid-90degTilt5degRes" },
    @{ Row = 22; A = 20; B = "HexGrid-90degTilt10degRes" },
    @{ Row = 23; A = 21; B = "HexGrid-90degTilt15degRes" }
)

foreach ($r in $newRows) {
    $rowNum = $r.Row

    # Copy formatting from the row above so the new row matches the existing
    # table styling (bold/bordered/centered column A, plain data cells).
    $ws.Range("A" + ($rowNum - 1) + ":T" + ($rowNum - 1)).Copy()
    $ws.Range("A" + $rowNum + ":T" + $rowNum).PasteSpecial(-4122)

    $ws.Range("A" + $rowNum).Value = $r.A
    $ws.Range("B" + $rowNum).Value = $r.B

    for ($col = 3; $col -le 20; $col++) {
        $ws.Cells.Item($rowNum, $col).Value = 1
    }
}

$excel.CutCopyMode = 0
